$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data block between rows 2-5 (La Ligua / 44890) and rows 6-9
# (Provincia de Limarí / 44908) for columns D, M, N, O, P, R, S.
$pairs = @(
    @{ a = 2; b = 6 },
    @{ a = 3; b = 7 },
    @{ a = 4; b = 8 },
    @{ a = 5; b = 9 }
)

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($pair in $pairs) {
    $rowA = $pair.a
    $rowB = $pair.b
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valueA = $rangeA.Value2
        $valueB = $rangeB.Value2
        $rangeA.Value2 = $valueB
        $rangeB.Value2 = $valueA
    }
}
